# Scheduled runner update: refresh market-price-derived profit figures
# across the per-job profit sheets (currentAveragePrice / derived
# totals in columns H-N). Only the affected cells are touched; all
# other data (item names, levels, leve gil, etc.) is left untouched.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 8658.5
$ws.Range("I18").Value = 9790.200000000001
$ws.Range("K18").Value = 9790.200000000001
$ws.Range("M18").Value = -9506.200000000001

$ws.Range("H41").Value = 45771.227
$ws.Range("I41").Value = 86.77778000000001
$ws.Range("J41").Value = 77398.92
$ws.Range("K41").Value = 86.77778000000001
$ws.Range("L41").Value = 77398.92
$ws.Range("M41").Value = 353.22222
$ws.Range("N41").Value = -78278.92

$ws.Range("H62").Value = 15671.695
$ws.Range("I62").Value = 12872.45
$ws.Range("J62").Value = 34333.332
$ws.Range("K62").Value = 12872.45
$ws.Range("L62").Value = 34333.332
$ws.Range("M62").Value = -12248.45
$ws.Range("N62").Value = -35581.332

$ws.Range("H65").Value = 15671.695
$ws.Range("I65").Value = 12872.45
$ws.Range("J65").Value = 34333.332
$ws.Range("K65").Value = 64362.25
$ws.Range("L65").Value = 171666.66
$ws.Range("M65").Value = -61242.25
$ws.Range("N65").Value = -177906.66

$ws.Range("H69").Value = 8006.222
$ws.Range("I69").Value = 3999.3333
$ws.Range("J69").Value = 10009.667
$ws.Range("K69").Value = 11997.9999
$ws.Range("L69").Value = 30029.001
$ws.Range("M69").Value = -11123.9999
$ws.Range("N69").Value = -31777.001

$ws.Range("H72").Value = 8006.222
$ws.Range("I72").Value = 3999.3333
$ws.Range("J72").Value = 10009.667
$ws.Range("K72").Value = 35993.9997
$ws.Range("L72").Value = 90087.003
$ws.Range("M72").Value = -31625.9997
$ws.Range("N72").Value = -98823.003

$ws.Range("H137").Value = 1128889.8
$ws.Range("I137").Value = 993408.5600000001
$ws.Range("K137").Value = 2980225.68
$ws.Range("M137").Value = -2977675.68


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4572349.5
$ws.Range("I32").Value = 5133115.5
$ws.Range("K32").Value = 5133115.5
$ws.Range("M32").Value = -5132828.5

$ws.Range("H81").Value = 80181
$ws.Range("J81").Value = 80181
$ws.Range("L81").Value = 80181
$ws.Range("N81").Value = -82177

$ws.Range("H84").Value = 80181
$ws.Range("J84").Value = 80181
$ws.Range("L84").Value = 240543
$ws.Range("N84").Value = -250527

$ws.Range("H102").Value = 4000.7646
$ws.Range("I102").Value = 3627.75
$ws.Range("K102").Value = 3627.75
$ws.Range("M102").Value = -2005.75

$ws.Range("H132").Value = 254747.39
$ws.Range("I132").Value = 418988.66
$ws.Range("K132").Value = 1256965.98
$ws.Range("M132").Value = -1254435.98


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1273.8889
$ws.Range("I86").Value = 1273.8889
$ws.Range("K86").Value = 1273.8889
$ws.Range("M86").Value = -150.8888999999999

$ws.Range("H89").Value = 1273.8889
$ws.Range("I89").Value = 1273.8889
$ws.Range("K89").Value = 6369.4445
$ws.Range("M89").Value = -753.4444999999996

$ws.Range("H107").Value = 1960.625
$ws.Range("I107").Value = 2106.2
$ws.Range("J107").Value = 1232.75
$ws.Range("K107").Value = 2106.2
$ws.Range("L107").Value = 1232.75
$ws.Range("M107").Value = -186.1999999999998
$ws.Range("N107").Value = -5072.75


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7413.9375
$ws.Range("I31").Value = 1844.7693
$ws.Range("J31").Value = 9482.485000000001
$ws.Range("K31").Value = 1844.7693
$ws.Range("L31").Value = 9482.485000000001
$ws.Range("M31").Value = -1549.7693
$ws.Range("N31").Value = -10072.485

$ws.Range("H34").Value = 7413.9375
$ws.Range("I34").Value = 1844.7693
$ws.Range("J34").Value = 9482.485000000001
$ws.Range("K34").Value = 1844.7693
$ws.Range("L34").Value = 9482.485000000001
$ws.Range("M34").Value = -1642.7693
$ws.Range("N34").Value = -9886.485000000001

$ws.Range("H39").Value = 3964.6667
$ws.Range("I39").Value = 3964.6667
$ws.Range("K39").Value = 3964.6667
$ws.Range("M39").Value = -3573.6667

$ws.Range("H49").Value = 3964.6667
$ws.Range("I49").Value = 3964.6667
$ws.Range("K49").Value = 3964.6667
$ws.Range("M49").Value = -3782.6667

$ws.Range("H99").Value = 8064
$ws.Range("I99").Value = 8648.111000000001
$ws.Range("J99").Value = 6749.75
$ws.Range("K99").Value = 8648.111000000001
$ws.Range("L99").Value = 6749.75
$ws.Range("M99").Value = -7150.111000000001
$ws.Range("N99").Value = -9745.75

$ws.Range("H122").Value = 993.2857
$ws.Range("I122").Value = 877.35297
$ws.Range("K122").Value = 2632.05891
$ws.Range("M122").Value = -182.0589100000002

$ws.Range("H126").Value = 8064
$ws.Range("I126").Value = 8648.111000000001
$ws.Range("J126").Value = 6749.75
$ws.Range("K126").Value = 25944.333
$ws.Range("L126").Value = 20249.25
$ws.Range("M126").Value = -23474.333
$ws.Range("N126").Value = -25189.25

$ws.Range("H134").Value = 12253.594
$ws.Range("I134").Value = 12622.549
$ws.Range("K134").Value = 37867.647
$ws.Range("M134").Value = -35332.647


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 24677.5
$ws.Range("I120").Value = 14336.143
$ws.Range("K120").Value = 43008.429
$ws.Range("M120").Value = -38170.429


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 203308.31
$ws.Range("I80").Value = 266350.78
$ws.Range("K80").Value = 266350.78
$ws.Range("M80").Value = -265352.78

$ws.Range("H83").Value = 203308.31
$ws.Range("I83").Value = 266350.78
$ws.Range("K83").Value = 1331753.9
$ws.Range("M83").Value = -1326761.9

$ws.Range("H102").Value = 3876.05
$ws.Range("I102").Value = 3418.111
$ws.Range("K102").Value = 3418.111
$ws.Range("M102").Value = -1796.111

$ws.Range("H136").Value = 47177
$ws.Range("J136").Value = 47177
$ws.Range("L136").Value = 141531
$ws.Range("N136").Value = -146631


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 55591.156
$ws.Range("I22").Value = 200797.6
$ws.Range("J22").Value = 3731.7144
$ws.Range("K22").Value = 200797.6
$ws.Range("L22").Value = 3731.7144
$ws.Range("M22").Value = -200502.6
$ws.Range("N22").Value = -4321.7144

$ws.Range("H27").Value = 55591.156
$ws.Range("I27").Value = 200797.6
$ws.Range("J27").Value = 3731.7144
$ws.Range("K27").Value = 200797.6
$ws.Range("L27").Value = 3731.7144
$ws.Range("M27").Value = -200690.6
$ws.Range("N27").Value = -3945.7144

$ws.Range("H61").Value = 16414.562
$ws.Range("I61").Value = 18488.143
$ws.Range("J61").Value = 1899.5
$ws.Range("K61").Value = 18488.143
$ws.Range("L61").Value = 1899.5
$ws.Range("M61").Value = -18286.143
$ws.Range("N61").Value = -2303.5

$ws.Range("H68").Value = 5772.35
$ws.Range("I68").Value = 4241.1665
$ws.Range("J68").Value = 6428.5713
$ws.Range("K68").Value = 4241.1665
$ws.Range("L68").Value = 6428.5713
$ws.Range("M68").Value = -3492.1665
$ws.Range("N68").Value = -7926.5713

$ws.Range("H71").Value = 5772.35
$ws.Range("I71").Value = 4241.1665
$ws.Range("J71").Value = 6428.5713
$ws.Range("K71").Value = 21205.8325
$ws.Range("L71").Value = 32142.8565
$ws.Range("M71").Value = -17461.8325
$ws.Range("N71").Value = -39630.85649999999

$ws.Range("H113").Value = 16414.562
$ws.Range("I113").Value = 18488.143
$ws.Range("J113").Value = 1899.5
$ws.Range("K113").Value = 18488.143
$ws.Range("L113").Value = 1899.5
$ws.Range("M113").Value = -16318.143
$ws.Range("N113").Value = -6239.5


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1893.6666
$ws.Range("I113").Value = 1507.1666
$ws.Range("J113").Value = 2666.6667
$ws.Range("K113").Value = 4521.4998
$ws.Range("L113").Value = 8000.000100000001
$ws.Range("M113").Value = -2351.4998

